$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 468.16666
$ws.Range("I9").Value = 401.8
$ws.Range("K9").Value = 401.8
$ws.Range("M9").Value = -232.8
$ws.Range("H17").Value = 2358.6
$ws.Range("J17").Value = 2358.6
$ws.Range("L17").Value = 7075.799999999999
$ws.Range("N17").Value = -7411.799999999999
$ws.Range("H43").Value = 2071.3635
$ws.Range("I43").Value = 2479.8
$ws.Range("J43").Value = 1731
$ws.Range("K43").Value = 2479.8
$ws.Range("L43").Value = 1731
$ws.Range("M43").Value = -2410.8
$ws.Range("N43").Value = -1869
$ws.Range("H76").Value = 6761.533
$ws.Range("I76").Value = 5942.5557
$ws.Range("K76").Value = 5942.5557
$ws.Range("M76").Value = -5627.5557
$ws.Range("H79").Value = 6761.533
$ws.Range("I79").Value = 5942.5557
$ws.Range("K79").Value = 5942.5557
$ws.Range("M79").Value = -4850.5557
$ws.Range("H80").Value = 4167276
$ws.Range("I80").Value = 8928837
$ws.Range("J80").Value = 910
$ws.Range("K80").Value = 26786511
$ws.Range("L80").Value = 2730
$ws.Range("M80").Value = -26785513
$ws.Range("N80").Value = -4726
$ws.Range("H83").Value = 4167276
$ws.Range("I83").Value = 8928837
$ws.Range("J83").Value = 910
$ws.Range("K83").Value = 80359533
$ws.Range("L83").Value = 8190
$ws.Range("M83").Value = -80354541
$ws.Range("N83").Value = -18174
$ws.Range("H87").Value = 76314
$ws.Range("J87").Value = 91931.75
$ws.Range("L87").Value = 91931.75
$ws.Range("N87").Value = -94427.75
$ws.Range("H90").Value = 76314
$ws.Range("J90").Value = 91931.75
$ws.Range("L90").Value = 275795.25
$ws.Range("N90").Value = -288275.25
$ws.Range("H112").Value = 1473.5682
$ws.Range("J112").Value = 1503.0238
$ws.Range("L112").Value = 4509.0714
$ws.Range("N112").Value = -6725.0714
$ws.Range("H133").Value = 80779
$ws.Range("J133").Value = 80779
$ws.Range("L133").Value = 80779
$ws.Range("N133").Value = -90899
$ws.Range("H137").Value = 1032593.2
$ws.Range("I137").Value = 731689.3
$ws.Range("K137").Value = 2195067.9
$ws.Range("M137").Value = -2192517.9
$ws.Range("H138").Value = 2726.0544
$ws.Range("I138").Value = 2008.1818
$ws.Range("J138").Value = 3204.6365
$ws.Range("K138").Value = 6024.5454
$ws.Range("L138").Value = 9613.9095
$ws.Range("M138").Value = -884.5454
$ws.Range("N138").Value = -19893.9095
$ws.Range("H139").Value = 95000
$ws.Range("J139").Value = 95000
$ws.Range("L139").Value = 95000
$ws.Range("N139").Value = -105280

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3136.5715
$ws.Range("J2").Value = 4650.4
$ws.Range("L2").Value = 4650.4
$ws.Range("N2").Value = -4876.4
$ws.Range("H116").Value = 3136.5715
$ws.Range("J116").Value = 4650.4
$ws.Range("L116").Value = 4650.4
$ws.Range("N116").Value = -9238.4
$ws.Range("H122").Value = 3734.5
$ws.Range("I122").Value = 3263.6428
$ws.Range("J122").Value = 4833.1665
$ws.Range("K122").Value = 9790.928400000001
$ws.Range("L122").Value = 14499.4995
$ws.Range("M122").Value = -7340.928400000001
$ws.Range("N122").Value = -19399.4995

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3136.5715
$ws.Range("J3").Value = 4650.4
$ws.Range("L3").Value = 4650.4
$ws.Range("N3").Value = -4878.4
$ws.Range("H20").Value = 3115.375
$ws.Range("I20").Value = 4521.3335
$ws.Range("J20").Value = 2271.8
$ws.Range("K20").Value = 4521.3335
$ws.Range("L20").Value = 2271.8
$ws.Range("M20").Value = -4274.3335
$ws.Range("N20").Value = -2765.8
$ws.Range("H96").Value = 10229
$ws.Range("I96").Value = 10229
$ws.Range("K96").Value = 10229
$ws.Range("M96").Value = -7483
$ws.Range("H99").Value = 30387.578

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 65332.145
$ws.Range("J68").Value = 68054.164
$ws.Range("L68").Value = 68054.164
$ws.Range("N68").Value = -69552.164
$ws.Range("H71").Value = 65332.145
$ws.Range("J71").Value = 68054.164
$ws.Range("L71").Value = 204162.492
$ws.Range("N71").Value = -211650.492

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 217426
$ws.Range("I2").Value = 333361.66
$ws.Range("J2").Value = 46.625
$ws.Range("K2").Value = 2000169.96
$ws.Range("L2").Value = 279.75
$ws.Range("M2").Value = -2000056.96
$ws.Range("N2").Value = -505.75
$ws.Range("H3").Value = 5853.4375
$ws.Range("I3").Value = 1432.8334
$ws.Range("J3").Value = 19115.25
$ws.Range("K3").Value = 4298.5002
$ws.Range("L3").Value = 57345.75
$ws.Range("M3").Value = -4186.5002
$ws.Range("N3").Value = -57569.75
$ws.Range("H37").Value = 85555.75
$ws.Range("J37").Value = 85555.75
$ws.Range("L37").Value = 256667.25
$ws.Range("N37").Value = -256891.25
$ws.Range("H50").Value = 730.75
$ws.Range("I50").Value = 269.2
$ws.Range("J50").Value = 1500
$ws.Range("K50").Value = 807.5999999999999
$ws.Range("L50").Value = 4500
$ws.Range("M50").Value = -326.5999999999999
$ws.Range("N50").Value = -5462
$ws.Range("H53").Value = 730.75
$ws.Range("I53").Value = 269.2
$ws.Range("J53").Value = 1500
$ws.Range("K53").Value = 807.5999999999999
$ws.Range("L53").Value = 4500
$ws.Range("M53").Value = -326.5999999999999
$ws.Range("N53").Value = -5462
$ws.Range("H56").Value = 7825.273
$ws.Range("I56").Value = 7825.273
$ws.Range("K56").Value = 7825.273
$ws.Range("M56").Value = -7295.273
$ws.Range("H68").Value = 9442.5
$ws.Range("I68").Value = 2423.75
$ws.Range("J68").Value = 12951.875
$ws.Range("K68").Value = 7271.25
$ws.Range("L68").Value = 38855.625
$ws.Range("M68").Value = -6460.25
$ws.Range("N68").Value = -40477.625
$ws.Range("H71").Value = 9442.5
$ws.Range("I71").Value = 2423.75
$ws.Range("J71").Value = 12951.875
$ws.Range("K71").Value = 21813.75
$ws.Range("L71").Value = 116566.875
$ws.Range("M71").Value = -17757.75
$ws.Range("N71").Value = -124678.875
$ws.Range("H80").Value = 4632
$ws.Range("I80").Value = 4000
$ws.Range("J80").Value = 4790
$ws.Range("K80").Value = 12000
$ws.Range("L80").Value = 14370
$ws.Range("M80").Value = -11064
$ws.Range("N80").Value = -16242
$ws.Range("H83").Value = 4632
$ws.Range("I83").Value = 4000
$ws.Range("J83").Value = 4790
$ws.Range("K83").Value = 36000
$ws.Range("L83").Value = 43110
$ws.Range("M83").Value = -31320
$ws.Range("N83").Value = -52470
$ws.Range("H113").Value = 1688.4286
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1688.4286
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 5065.2858
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -9405.2858
$ws.Range("H117").Value = 5741.75
$ws.Range("J117").Value = 8064
$ws.Range("L117").Value = 24192
$ws.Range("N117").Value = -31076

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18148.2
$ws.Range("J70").Value = 9998.5
$ws.Range("L70").Value = 9998.5
$ws.Range("N70").Value = -10538.5
$ws.Range("H73").Value = 18148.2
$ws.Range("J73").Value = 9998.5
$ws.Range("L73").Value = 9998.5
$ws.Range("N73").Value = -11870.5
$ws.Range("H80").Value = 6214
$ws.Range("J80").Value = 6487.375
$ws.Range("L80").Value = 6487.375
$ws.Range("N80").Value = -8483.375
$ws.Range("H83").Value = 6214
$ws.Range("J83").Value = 6487.375
$ws.Range("L83").Value = 32436.875
$ws.Range("N83").Value = -42420.875
$ws.Range("H122").Value = 3491.8262
$ws.Range("I122").Value = 3296.6667
$ws.Range("J122").Value = 3857.75
$ws.Range("K122").Value = 9890.000100000001
$ws.Range("L122").Value = 11573.25
$ws.Range("M122").Value = -7440.000100000001
$ws.Range("N122").Value = -16473.25

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H122").Value = 3353.0513
$ws.Range("I122").Value = 3186.074
$ws.Range("J122").Value = 3728.75
$ws.Range("K122").Value = 9558.222
$ws.Range("L122").Value = 11186.25
$ws.Range("M122").Value = -7108.222
$ws.Range("N122").Value = -16086.25
$ws.Range("H136").Value = 6485.0557
$ws.Range("I136").Value = 5840.25
$ws.Range("J136").Value = 7774.6665
$ws.Range("K136").Value = 17520.75
$ws.Range("L136").Value = 23323.9995
$ws.Range("M136").Value = -14970.75
$ws.Range("N136").Value = -28423.9995

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1783.0588
$ws.Range("I122").Value = 1343.8182
$ws.Range("J122").Value = 2588.3333
$ws.Range("K122").Value = 4031.4546
$ws.Range("L122").Value = 7764.999899999999
$ws.Range("M122").Value = -1581.4546
$ws.Range("N122").Value = -12664.9999
